# Automatische test-sync: 2025-07-23 22:22:50
#
# Appends a new test-mail row (#4) to the "Logs" sheet and a matching
# roll-up row to the "Dashboard" sheet, then widens the chart series /
# conditional-formatting ranges so they keep covering the new data.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Logs!A14:J14 --------------------------------------------------------
$logs.Range("A14").Value = "Wil je deze factuur even printen en klaarleggen voor maandag?"
$logs.Range("B14").Value = "mailmind.test@zohomail.eu"
$logs.Range("C14").Value = "Testmail #4: Wil je deze factuur even printen en klaarleggen voor maandag?"
$logs.Range("D14").Value = "Factuur / Administratie"
$logs.Range("E14").Value = "Geachte heer/mevrouw,`nHelaas kan ik u niet helpen met het printen en klaarleggen van de factuur. Ik ben een e-mailassistent en kan geen fysieke taken uitvoeren. Ik raad u aan contact op te nemen met de relevante afdeling binnen uw bedrijf om dit verzoek in behandeling te nemen.`nMet vriendelijke groet,`n[Je Naam]`nE-mailassistent"
$logs.Range("F14").Value = "2025-07-23 22:22:32"
$logs.Range("G14").Value = "Ja"
$logs.Range("H14").Value = "Ja"
$logs.Range("I14").Value = "Nee"
$logs.Range("J14").Value = "Nee"

# --- Dashboard!A5:B5 -------------------------------------------------------
$dash.Range("A5").Value = "Factuur / Administratie"
$dash.Range("B5").Value = 1

# --- Widen the conditional-formatting rules on "Logs" from row 13 to 14 ---
$ranges = @("D2:D13", "G2:G13", "H2:H13", "I2:I13", "J2:J13")
foreach ($rng in $ranges) {
    $col = $rng.Substring(0, 1)
    $newRange = $logs.Range($col + "2:" + $col + "14")
    $fcs = $logs.Range($rng).FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Widen the "Dashboard" chart's category/value series to row 5 ---------
$co = $dash.ChartObjects().Item(1)
$chart = $co.Chart
$ser = $chart.SeriesCollection().Item(1)
$ser.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$5,'Dashboard'!`$B`$2:`$B`$5,1)"

$wb.Save()
